$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '66.948.22'
Set-TextCell 2 5 '  -5.02%  '
Set-TextCell 3 4 '3.219.55'
Set-TextCell 3 5 '  -8.61%  '
Set-TextCell 4 4 '0.999'
Set-TextCell 4 5 '  -0.11%  '
Set-TextCell 5 4 '592.19'
Set-TextCell 5 5 '  -2.47%  '
Set-TextCell 6 4 '151.90'
Set-TextCell 6 5 '  -12.88%  '
Set-TextCell 7 5 '  -0.02%  '
Set-TextCell 8 4 '3.209.27'
Set-TextCell 8 5 '  -8.79%  '
Set-TextCell 9 5 '  -11.51%  '
Set-TextCell 10 5 '  -12.67%  '
Set-TextCell 11 4 '6.44'
Set-TextCell 11 5 '  -10.91%  '
Set-TextCell 12 5 '  -15.99%  '
Set-TextCell 13 4 '38.92'
Set-TextCell 13 5 '  -16.04%  '
Set-TextCell 14 4 '0.0000243'
Set-TextCell 14 5 '  -12.58%  '
Set-TextCell 15 4 '3.749.84'
Set-TextCell 15 5 '  -8.45%  '
Set-TextCell 16 4 '66.974.69'
Set-TextCell 16 5 '  -5.02%  '
Set-TextCell 17 4 '3.230.56'
Set-TextCell 17 5 '  -8.76%  '
Set-TextCell 18 4 '0.114'
Set-TextCell 18 5 '  -4.59%  '
Set-TextCell 19 4 '528.88'
Set-TextCell 19 5 '  -13.92%  '
Set-TextCell 20 4 '7.09'
Set-TextCell 20 5 '  -14.86%  '
Set-TextCell 21 4 '14.85'
Set-TextCell 21 5 '  -15.23%  '
Set-TextCell 22 5 '  -14.11%  '
Set-TextCell 23 5 '  -12.46%  '
Set-TextCell 24 4 '13.83'
Set-TextCell 24 5 '  -11.67%  '
Set-TextCell 25 4 '85.48'
Set-TextCell 25 5 '  -14.04%  '
Set-TextCell 26 4 '0.998'
Set-TextCell 26 5 '  -0.20%  '
Set-TextCell 27 5 '  -14.69%  '
Set-TextCell 28 2 'RenderToken'
Set-TextCell 28 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 28 4 '8.14'
Set-TextCell 28 5 '  -10.23%  '
Set-TextCell 29 2 'ImmutableX'
Set-TextCell 29 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 29 4 '2.17'
Set-TextCell 29 5 '  -15.63%  '
Set-TextCell 30 4 '29.10'
Set-TextCell 30 5 '  -15.19%  '
Set-TextCell 31 4 '2.68'
Set-TextCell 31 5 '  -10.05%  '
Set-TextCell 32 5 '  -10.05%  '
Set-TextCell 33 4 '545.19'
Set-TextCell 33 5 '  -16.29%  '
Set-TextCell 34 4 '5.73'
Set-TextCell 34 5 '  -16.26%  '
Set-TextCell 35 4 '6.45'
Set-TextCell 35 5 '  -20.27%  '
Set-TextCell 36 5 '  +0.30%  '
Set-TextCell 37 4 '53.63'
Set-TextCell 37 5 '  -5.63%  '
Set-TextCell 38 4 '0.0425'
Set-TextCell 38 5 '  -11.12%  '
Set-TextCell 39 4 '0.0859'
Set-TextCell 39 5 '  -14.04%  '
Set-TextCell 40 5 '  -13.90%  '
Set-TextCell 41 5 '  -12.55%  '
Set-TextCell 42 4 '2.916.69'
Set-TextCell 42 5 '  -13.49%  '
Set-TextCell 43 4 '2.68'
Set-TextCell 43 5 '  -24.77%  '
Set-TextCell 44 5 '  -15.18%  '
Set-TextCell 45 4 '0.0₃0580'
Set-TextCell 45 5 '  -22.70%  '
Set-TextCell 46 4 '2.44'
Set-TextCell 46 5 '  -16.03%  '
Set-TextCell 47 2 'Fetch.AI'
Set-TextCell 47 3 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell 47 4 '2.16'
Set-TextCell 47 5 '  -15.64%  '
Set-TextCell 48 2 'InjectiveProtocol'
Set-TextCell 48 3 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell 48 4 '26.32'
Set-TextCell 48 5 '  -18.38%  '
Set-TextCell 50 4 '0.114'
Set-TextCell 50 5 '  -12.14%  '
Set-TextCell 51 4 '118.09'
Set-TextCell 51 5 '  -11.49%  '
